$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9319836497306824
$ws.Range("B1").Value = 3.162780284881592
$ws.Range("C1").Value = 4.20416259765625
$ws.Range("D1").Value = 3.024479389190674
$ws.Range("E1").Value = 1.369701504707336
